$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.869.49"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "3.760.01"
$ws.Range("E3").Value = "  -2.78%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'595.44"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").Value = "'167.23"
$ws.Range("E6").Value = "  -3.43%  "
$ws.Range("D7").Value = "3.760.79"
$ws.Range("E7").Value = "  -2.72%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.521"
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("D10").Value = "'0.162"
$ws.Range("E10").Value = "  -5.09%  "
$ws.Range("D11").Value = "'6.46"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("D12").Value = "'0.448"
$ws.Range("E12").Value = "  -3.38%  "
$ws.Range("D13").Value = "'0.0000271"
$ws.Range("E13").Value = "  -7.60%  "
$ws.Range("D14").Value = "'36.27"
$ws.Range("E14").Value = "  -3.04%  "
$ws.Range("D15").Value = "4.390.60"
$ws.Range("E15").Value = "  -2.79%  "
$ws.Range("D16").Value = "3.765.66"
$ws.Range("E16").Value = "  -2.88%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "67.804.57"
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'18.46"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").Value = "'7.08"
$ws.Range("E19").Value = "  -5.72%  "
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "'10.49"
$ws.Range("E21").Value = "  -4.51%  "
$ws.Range("D22").Value = "'466.02"
$ws.Range("E22").Value = "  -1.52%  "
$ws.Range("D23").Value = "'0.709"
$ws.Range("E23").Value = "  -3.63%  "
$ws.Range("D24").Value = "'83.41"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("D25").Value = "'0.0000144"
$ws.Range("E25").Value = "  -11.93%  "
$ws.Range("D26").Value = "'2.21"
$ws.Range("E26").Value = "  -3.92%  "
$ws.Range("D27").Value = "'12.07"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("D28").Value = "'10.30"
$ws.Range("E28").Value = "  -2.58%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'2.90"
$ws.Range("E30").Value = "  -1.95%  "
$ws.Range("D31").Value = "3.907.69"
$ws.Range("E31").Value = "  -2.74%  "
$ws.Range("D32").Value = "'7.47"
$ws.Range("E32").Value = "  -4.47%  "
$ws.Range("D33").Value = "'30.16"
$ws.Range("E33").Value = "  -3.78%  "
$ws.Range("D34").Value = "'2.20"
$ws.Range("E34").Value = "  -5.89%  "
$ws.Range("D35").Value = "'9.11"
$ws.Range("E35").Value = "  -4.12%  "
$ws.Range("D36").Value = "3.708.49"
$ws.Range("E36").Value = "  -3.25%  "
$ws.Range("E37").Value = "  -3.01%  "
$ws.Range("D38").Value = "'3.61"
$ws.Range("E38").Value = "  -8.78%  "
$ws.Range("D39").Value = "'0.138"
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("D41").Value = "'5.80"
$ws.Range("E41").Value = "  -3.91%  "
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.308"
$ws.Range("E44").Value = "  -4.71%  "
$ws.Range("D45").Value = "'8.58"
$ws.Range("E45").Value = "  -2.72%  "
$ws.Range("D46").Value = "'1.92"
$ws.Range("E46").Value = "  -4.22%  "
$ws.Range("D47").Value = "'399.01"
$ws.Range("E47").Value = "  -5.93%  "
$ws.Range("D48").Value = "'44.85"
$ws.Range("E48").Value = "  -3.95%  "
$ws.Range("D49").Value = "'143.70"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("D50").Value = "'39.05"
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("D51").Value = "'0.0349"
$ws.Range("E51").Value = "  -3.80%  "
